# Deploy 4ffca36 from branch develop
# Updates the "7.1.2" indicator sheet with the newly published 2017 data
# point, extends the Tabelle2 data table + chart series accordingly, and
# refreshes the WHO data-source citation.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")

# ---------------------------------------------------------------------
# 1) Tabelle1: add the 2017 row to the visible data table (row 12),
#    copying the formatting used by the previous years' rows.
# ---------------------------------------------------------------------
$ws1.Range("A11:B11").Copy()
$ws1.Range("A12:B12").PasteSpecial(-4122)
$ws1.Range("A12").Value = 2017
$ws1.Range("B12").Value = "> 95"

# ---------------------------------------------------------------------
# 2) Tabelle1: refresh the WHO data-source citation with the new date.
# ---------------------------------------------------------------------
$ws1.Range("B35").Value = "Weltgesundheitsorganisation (WHO), Stand 11.04.2019"

# ---------------------------------------------------------------------
# 3) Tabelle2 (chart source data): add the matching 2017 row.
# ---------------------------------------------------------------------
$ws2.Range("A7:C7").Copy()
$ws2.Range("A8:C8").PasteSpecial(-4122)
$ws2.Range("A8").Value = 2017
$ws2.Range("B8").Value = 95
$ws2.Range("C8").Value = 5

# ---------------------------------------------------------------------
# 4) Chart: extend both series to the new last row of Tabelle2.
# ---------------------------------------------------------------------
$co = $ws1.ChartObjects().Item(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection().Item(1)
$s1.Formula = "=SERIES(,Tabelle2!`$A`$1:`$A`$8,Tabelle2!`$B`$1:`$B`$8,1)"
$s2 = $chart.SeriesCollection().Item(2)
$s2.Formula = "=SERIES(,Tabelle2!`$A`$1:`$A`$8,Tabelle2!`$C`$1:`$C`$8,2)"

# ---------------------------------------------------------------------
# 5) Restore view state (selection / scroll position) as left by the
#    author after making the edit.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("C12").Select()

$ws1.Activate()
$ws1.Range("J24").Select()
